$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2251
$ws.Range("F3").Value = 276
$ws.Range("F4").Value = 165
$ws.Range("F5").Value = 163
$ws.Range("F6").Value = 305
$ws.Range("G6").Value = 54
$ws.Range("F8").Value = 677
$ws.Range("F9").Value = 498
$ws.Range("F10").Value = 612
$ws.Range("F11").Value = 364
$ws.Range("F12").Value = 60
$ws.Range("F13").Value = 343
$ws.Range("F14").Value = 950
$ws.Range("F15").Value = 220
$ws.Range("F17").Value = 94
$ws.Range("G17").Value = 60
$ws.Range("F19").Value = 6
$ws.Range("F21").Value = 218
$ws.Range("F22").Value = 86

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 78
$ws.Range("F6").Value = 166
$ws.Range("F7").Value = 200
$ws.Range("F8").Value = 2458
$ws.Range("F10").Value = 18
$ws.Range("F16").Value = 2296

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 247
$ws.Range("F3").Value = 9
$ws.Range("F4").Value = 335
$ws.Range("F5").Value = 162

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 247
$ws.Range("F3").Value = 9
$ws.Range("F6").Value = 2251
$ws.Range("F7").Value = 335
$ws.Range("F8").Value = 276
$ws.Range("F9").Value = 165
$ws.Range("F10").Value = 163
$ws.Range("F11").Value = 305
$ws.Range("G11").Value = 54
$ws.Range("F14").Value = 78
$ws.Range("F15").Value = 166
$ws.Range("F16").Value = 162
$ws.Range("F17").Value = 677
$ws.Range("F18").Value = 498
$ws.Range("F19").Value = 612
$ws.Range("F20").Value = 364
$ws.Range("F21").Value = 60
$ws.Range("F22").Value = 343
$ws.Range("F23").Value = 950
$ws.Range("F24").Value = 200
$ws.Range("F25").Value = 2458
$ws.Range("F27").Value = 18
$ws.Range("F31").Value = 220
$ws.Range("F33").Value = 94
$ws.Range("G33").Value = 60
$ws.Range("F37").Value = 6
$ws.Range("F39").Value = 218
$ws.Range("F40").Value = 86
$ws.Range("F41").Value = 2296
